$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.461.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.532"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.79%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.256"
$ws.Range("D9").Style = "Normal"

$ws.Range("E10").Value = "  -0.83%  "

$ws.Range("E11").Value = "  +1.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.868.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.633.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.577"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.17%  "

$ws.Range("E15").Value = "  -1.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.480.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.57%  "

$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.55%  "

$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.03%  "

$ws.Range("E23").Value = "  +5.42%  "

$ws.Range("E24").Value = "  -3.09%  "

$ws.Range("E25").Value = "  +2.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.69%  "

$ws.Range("E27").Value = "  +1.92%  "

$ws.Range("E28").Value = "  +0.01%  "

$ws.Range("E29").Value = "  -3.14%  "

$ws.Range("E30").Value = "  -0.54%  "

$ws.Range("E31").Value = "  -1.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.415.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.01%  "

$ws.Range("E35").Value = "  +3.05%  "

$ws.Range("E36").Value = "  -2.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.571"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0168"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.96%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.921"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +18.08%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.873"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.87%  "

$ws.Range("E41").Value = "  -1.01%  "

$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("E43").Value = "  +1.83%  "

$ws.Range("E44").Value = "  +0.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.778.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.68%  "

$ws.Range("E47").Value = "  -3.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.49%  "

$ws.Range("E49").Value = "  +1.19%  "

$ws.Range("E50").Value = "  -1.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.61%  "
